$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new values in column E (Code Review 2) for rows 2-6
$ws.Range("E2").Value = 25
$ws.Range("E3").Value = 25
$ws.Range("E4").Value = 25
$ws.Range("E5").Value = 25
$ws.Range("E6").Value = 100

# Update the selected cell on the sheet
$ws.Activate()
$ws.Range("B9").Select()
